# Rename sheets and relabel hosts/objects per the new naming convention.

$wb = $excel.ActiveWorkbook

# --- Sheet name changes ---
$wsFW  = $wb.Worksheets.Item(1)   # was "FortiGW"
$wsCP  = $wb.Worksheets.Item(2)   # was "CPMGMT"

$wsFW.Name = "Fortinet_1"
$wsCP.Name = "WANCPFW"

# --- "Fortinet_1" sheet (formerly FortiGW) cell relabels ---
$wsFW.Range("B2").Value = "Internal DB"
$wsFW.Range("B3").Value = "InternalDB_1"
$wsFW.Range("B4").Value = "InternalDB_1"
$wsFW.Range("A5").Value = "w10c"
$wsFW.Range("B5").Value = "InternalDB_1"

# --- "WANCPFW" sheet (formerly CPMGMT) cell relabels ---
$wsCP.Range("B3").Value = "DMZ"
$wsCP.Range("A4").Value = "ws12c"
$wsCP.Range("B4").Value = "InternalDB_1"
$wsCP.Range("A5").Value = "DMZ"
$wsCP.Range("B5").Value = "Fortinet_1"
$wsCP.Range("A6").Value = "webmain"
$wsCP.Range("B6").Value = "InternalDB_1"

# --- Refresh selection / active cell on each sheet (cosmetic UI state) ---
# Select the non-active sheet first so the FortiGW-derived sheet ends up as
# the active tab (it was the active tab before the edit too).
[void]$wsCP.Range("B5").Select()
[void]$wsFW.Range("G13").Select()
